# Handles float input without breaking stuff
#
# This marksheet worksheet previously left the "Student Ans" columns blank
# (the grading step that fills them in crashed on float scores before this
# fix). We now populate the Student-Ans columns with the student's actual
# answers, colour-code them against the Correct-Ans columns, fix up the
# summary row numbers/score, and drop the now-unused extra
# "page" of Student Ans/Correct Ans columns (G:H), since with the fix all
# answered/graded questions fit into the A:E block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-StyledCell {
    param($sheet, $addr, $styleName, $value)
    $cell = $sheet.Range($addr)
    if ($null -ne $value) {
        $cell.Value = $value
    }
    $cell.Style = $styleName
}

# ---------------------------------------------------------------------
# Summary block (rows 10-12): right/wrong/not-attempted counts, marking
# scheme and computed score now reflect the actual (float-safe) grading.
# ---------------------------------------------------------------------
Set-StyledCell $ws "A10" "mtitleStyle" $null
$ws.Range("B10").Value = 17
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = 28

Set-StyledCell $ws "A11" "mtitleStyle" $null
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

Set-StyledCell $ws "A12" "mtitleStyle" $null
$ws.Range("B12").Value = 68
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "66/112"

# ---------------------------------------------------------------------
# Drop the third Student Ans / Correct Ans page (columns G:H) - no longer
# needed now that answers are correctly distributed across A:E.
# ---------------------------------------------------------------------
$ws.Range("G15:H40").Clear()

# ---------------------------------------------------------------------
# Student Ans column (A) for every question row now shows what the
# student actually answered, colour coded:
#   correctStyle   (green) -> matches the Correct Ans column
#   incorrectStyle (red)   -> does not match the Correct Ans column
#   normalStyle    (black) -> left blank / not attempted
# ---------------------------------------------------------------------
Set-StyledCell $ws "A16" "incorrectStyle" "Option D"
Set-StyledCell $ws "A17" "normalStyle"    $null
Set-StyledCell $ws "A18" "correctStyle"   "Option B"
Set-StyledCell $ws "A19" "correctStyle"   "Option C"
Set-StyledCell $ws "A20" "normalStyle"    $null
Set-StyledCell $ws "A21" "correctStyle"   "Option C"
Set-StyledCell $ws "A22" "correctStyle"   "Option D"
Set-StyledCell $ws "A23" "normalStyle"    $null
Set-StyledCell $ws "A24" "normalStyle"    $null
Set-StyledCell $ws "A25" "correctStyle"   "Option A"
Set-StyledCell $ws "A26" "normalStyle"    $null
Set-StyledCell $ws "A27" "correctStyle"   "Option A"
Set-StyledCell $ws "A28" "normalStyle"    $null
Set-StyledCell $ws "A29" "normalStyle"    $null
Set-StyledCell $ws "A30" "correctStyle"   "Option B"
Set-StyledCell $ws "A31" "correctStyle"   "Option D"
Set-StyledCell $ws "A32" "correctStyle"   "Option C"
Set-StyledCell $ws "A33" "correctStyle"   "Option D"
Set-StyledCell $ws "A34" "incorrectStyle" "Option A"
Set-StyledCell $ws "A35" "correctStyle"   "Option D"
Set-StyledCell $ws "A36" "correctStyle"   "Option A"
Set-StyledCell $ws "A37" "normalStyle"    $null
Set-StyledCell $ws "A38" "correctStyle"   "Option A"
Set-StyledCell $ws "A39" "correctStyle"   "Option D"
Set-StyledCell $ws "A40" "normalStyle"    $null

# ---------------------------------------------------------------------
# Second Student Ans / Correct Ans page (D:E) only has 3 real questions
# (rows 16-18) once floats are handled correctly; the student answered all
# three correctly, and the trailing rows 19-40 are no longer used.
# ---------------------------------------------------------------------
Set-StyledCell $ws "D16" "correctStyle" "Option A"
Set-StyledCell $ws "D17" "correctStyle" "Option C"
Set-StyledCell $ws "D18" "correctStyle" "Option D"
$ws.Range("D19:E40").Clear()
